$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.047608429025024
$ws.Cells.Item(2, 4).Value = 1.044381863062638
$ws.Cells.Item(2, 5).Value = 1.054068460966527
$ws.Cells.Item(2, 6).Value = 1.062642762355889
$ws.Cells.Item(2, 9).Value = 1.034920297112269
$ws.Cells.Item(2, 10).Value = 1.052656590605688
$ws.Cells.Item(2, 11).Value = 1.047153071356502
$ws.Cells.Item(2, 12).Value = 1.056812690318958
$ws.Cells.Item(2, 13).Value = 1.065363555061269

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.049185259509894
$ws.Cells.Item(3, 4).Value = 1.045511126727465
$ws.Cells.Item(3, 5).Value = 1.055514147769661
$ws.Cells.Item(3, 6).Value = 1.064273750346214
$ws.Cells.Item(3, 9).Value = 1.035215496809748
$ws.Cells.Item(3, 10).Value = 1.053878881405362
$ws.Cells.Item(3, 11).Value = 1.048092731787553
$ws.Cells.Item(3, 12).Value = 1.058069929759843
$ws.Cells.Item(3, 13).Value = 1.066807348783989

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.050203793147684
$ws.Cells.Item(4, 4).Value = 1.046240063579248
$ws.Cells.Item(4, 5).Value = 1.056448248816177
$ws.Cells.Item(4, 6).Value = 1.065327974101146
$ws.Cells.Item(4, 9).Value = 1.035404403034406
$ws.Cells.Item(4, 10).Value = 1.054667647215959
$ws.Cells.Item(4, 11).Value = 1.048698405076114
$ws.Cells.Item(4, 12).Value = 1.058881578580478
$ws.Cells.Item(4, 13).Value = 1.067739965838105

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.050631567642782
$ws.Cells.Item(5, 4).Value = 1.046546089941671
$ws.Cells.Item(5, 5).Value = 1.05684062825651
$ws.Cells.Item(5, 6).Value = 1.065770907049535
$ws.Cells.Item(5, 9).Value = 1.035483316578883
$ws.Cells.Item(5, 10).Value = 1.054998739537171
$ws.Cells.Item(5, 11).Value = 1.048952472852035
$ws.Cells.Item(5, 12).Value = 1.059222355373945
$ws.Cells.Item(5, 13).Value = 1.068131659719771

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.050703368677852
$ws.Cells.Item(6, 4).Value = 1.046597448723947
$ws.Cells.Item(6, 5).Value = 1.056906492147342
$ws.Cells.Item(6, 6).Value = 1.065845262291988
$ws.Cells.Item(6, 9).Value = 1.03549653709225
$ws.Cells.Item(6, 10).Value = 1.055054301994996
$ws.Cells.Item(6, 11).Value = 1.048995099423451
$ws.Cells.Item(6, 12).Value = 1.05927954767925
$ws.Cells.Item(6, 13).Value = 1.068197404894331

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.050209510718585
$ws.Cells.Item(7, 4).Value = 1.046244154356771
$ws.Cells.Item(7, 5).Value = 1.056453493044084
$ws.Cells.Item(7, 6).Value = 1.065333893614205
$ws.Cells.Item(7, 9).Value = 1.03540545945365
$ws.Cells.Item(7, 10).Value = 1.054672073262721
$ws.Cells.Item(7, 11).Value = 1.04870180212539
$ws.Cells.Item(7, 12).Value = 1.058886133778266
$ws.Cells.Item(7, 13).Value = 1.067745201147536

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.048141700714155
$ws.Cells.Item(8, 4).Value = 1.044763873339296
$ws.Cells.Item(8, 5).Value = 1.054557322826591
$ws.Cells.Item(8, 6).Value = 1.063194202664956
$ws.Cells.Item(8, 9).Value = 1.035020498834129
$ws.Cells.Item(8, 10).Value = 1.053070116452544
$ws.Cells.Item(8, 11).Value = 1.047471123740856
$ws.Cells.Item(8, 12).Value = 1.057237971592182
$ws.Cells.Item(8, 13).Value = 1.065851830979738

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.044483893126027
$ws.Cells.Item(9, 4).Value = 1.04214160373407
$ws.Cells.Item(9, 5).Value = 1.051205302516635
$ws.Cells.Item(9, 6).Value = 1.059414704751302
$ws.Cells.Item(9, 9).Value = 1.034325917899638
$ws.Cells.Item(9, 10).Value = 1.050230572266004
$ws.Cells.Item(9, 11).Value = 1.045284276116232
$ws.Cells.Item(9, 12).Value = 1.054319063610507
$ws.Cells.Item(9, 13).Value = 1.062502746099586

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.042035309320701
$ws.Cells.Item(10, 4).Value = 1.040383773376753
$ws.Cells.Item(10, 5).Value = 1.048962920825392
$ws.Cells.Item(10, 6).Value = 1.056888366954754
$ws.Cells.Item(10, 9).Value = 1.033851825778358
$ws.Cells.Item(10, 10).Value = 1.048325875183271
$ws.Cells.Item(10, 11).Value = 1.043813782017492
$ws.Cells.Item(10, 12).Value = 1.052362834888256
$ws.Cells.Item(10, 13).Value = 1.060260957184982

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040972530728777
$ws.Cells.Item(11, 4).Value = 1.039620244514013
$ws.Cells.Item(11, 5).Value = 1.047990008384568
$ws.Cells.Item(11, 6).Value = 1.055792726526435
$ws.Cells.Item(11, 9).Value = 1.033643892218008
$ws.Cells.Item(11, 10).Value = 1.047498253127325
$ws.Cells.Item(11, 11).Value = 1.043173978901835
$ws.Cells.Item(11, 12).Value = 1.05151322477629
$ws.Cells.Item(11, 13).Value = 1.059287973336831

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.040577376525442
$ws.Cells.Item(12, 4).Value = 1.039336271838596
$ws.Cells.Item(12, 5).Value = 1.047628323429212
$ws.Cells.Item(12, 6).Value = 1.055385487652794
$ws.Cells.Item(12, 9).Value = 1.033566256014586
$ws.Cells.Item(12, 10).Value = 1.047190396723542
$ws.Cells.Item(12, 11).Value = 1.04293585979148
$ws.Cells.Item(12, 12).Value = 1.051197249847286
$ws.Cells.Item(12, 13).Value = 1.058926212176779

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.040662156318248
$ws.Cells.Item(13, 4).Value = 1.039397201528774
$ws.Cells.Item(13, 5).Value = 1.047705919898048
$ws.Cells.Item(13, 6).Value = 1.055472854133085
$ws.Cells.Item(13, 9).Value = 1.033582927404178
$ws.Cells.Item(13, 10).Value = 1.047256453051141
$ws.Cells.Item(13, 11).Value = 1.042986958421982
$ws.Cells.Item(13, 12).Value = 1.051265045418319
$ws.Cells.Item(13, 13).Value = 1.05900382724623

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.040939875200263
$ws.Cells.Item(14, 4).Value = 1.039596778711433
$ws.Cells.Item(14, 5).Value = 1.047960117608603
$ws.Cells.Item(14, 6).Value = 1.055759069589815
$ws.Cells.Item(14, 9).Value = 1.033637482968276
$ws.Cells.Item(14, 10).Value = 1.047472814662096
$ws.Cells.Item(14, 11).Value = 1.043154305482135
$ws.Cells.Item(14, 12).Value = 1.051487114235841
$ws.Cells.Item(14, 13).Value = 1.059258077277501

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.041110934824776
$ws.Cells.Item(15, 4).Value = 1.039719696348988
$ws.Cells.Item(15, 5).Value = 1.048116696908854
$ws.Cells.Item(15, 6).Value = 1.055935380375001
$ws.Cells.Item(15, 9).Value = 1.033671043316778
$ws.Cells.Item(15, 10).Value = 1.047606063530184
$ws.Cells.Item(15, 11).Value = 1.043257351349839
$ws.Cells.Item(15, 12).Value = 1.051623885968627
$ws.Cells.Item(15, 13).Value = 1.059414682219369

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.042105788246933
$ws.Cells.Item(16, 4).Value = 1.040434395603407
$ws.Cells.Item(16, 5).Value = 1.049027448000087
$ws.Cells.Item(16, 6).Value = 1.056961043845565
$ws.Cells.Item(16, 9).Value = 1.03386556962877
$ws.Cells.Item(16, 10).Value = 1.04838074040696
$ws.Cells.Item(16, 11).Value = 1.043856178411661
$ws.Cells.Item(16, 12).Value = 1.052419166237335
$ws.Cells.Item(16, 13).Value = 1.060325482153172

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.042729149265117
$ws.Cells.Item(17, 4).Value = 1.040882066338727
$ws.Cells.Item(17, 5).Value = 1.049598210293363
$ws.Cells.Item(17, 6).Value = 1.057603947419795
$ws.Cells.Item(17, 9).Value = 1.033986879992691
$ws.Cells.Item(17, 10).Value = 1.048865898742429
$ws.Cells.Item(17, 11).Value = 1.044230980388098
$ws.Cells.Item(17, 12).Value = 1.052917335619225
$ws.Cells.Item(17, 13).Value = 1.060896186972301

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.043092502292185
$ws.Cells.Item(18, 4).Value = 1.041142956261732
$ws.Cells.Item(18, 5).Value = 1.049930939134834
$ws.Cells.Item(18, 6).Value = 1.057978777033222
$ws.Cells.Item(18, 9).Value = 1.034057382900526
$ws.Cells.Item(18, 10).Value = 1.049148606524112
$ws.Cells.Item(18, 11).Value = 1.044449300329842
$ws.Cells.Item(18, 12).Value = 1.05320766393418
$ws.Cells.Item(18, 13).Value = 1.061228850661203

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.043216355428199
$ws.Cells.Item(19, 4).Value = 1.041231874377766
$ws.Cells.Item(19, 5).Value = 1.050044359614786
$ws.Cells.Item(19, 6).Value = 1.058106556558348
$ws.Cells.Item(19, 9).Value = 1.034081379329437
$ws.Cells.Item(19, 10).Value = 1.049244955871821
$ws.Cells.Item(19, 11).Value = 1.044523691828837
$ws.Cells.Item(19, 12).Value = 1.053306617016992
$ws.Cells.Item(19, 13).Value = 1.061342243572779

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.042662293735415
$ws.Cells.Item(20, 4).Value = 1.04083405922642
$ws.Cells.Item(20, 5).Value = 1.049536992327527
$ws.Cells.Item(20, 6).Value = 1.05753498711329
$ws.Cells.Item(20, 9).Value = 1.03397389097151
$ws.Cells.Item(20, 10).Value = 1.048813874545088
$ws.Cells.Item(20, 11).Value = 1.044190798317523
$ws.Cells.Item(20, 12).Value = 1.052863912211451
$ws.Cells.Item(20, 13).Value = 1.060834978447236

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.040858104795176
$ws.Cells.Item(21, 4).Value = 1.039538018311732
$ws.Cells.Item(21, 5).Value = 1.047885271131019
$ws.Cells.Item(21, 6).Value = 1.055674793791766
$ws.Cells.Item(21, 9).Value = 1.03362142878825
$ws.Cells.Item(21, 10).Value = 1.047409113796602
$ws.Cells.Item(21, 11).Value = 1.043105038909277
$ws.Cells.Item(21, 12).Value = 1.051421731393312
$ws.Cells.Item(21, 13).Value = 1.059183216767349

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.039721472141338
$ws.Cells.Item(22, 4).Value = 1.038721036102713
$ws.Cells.Item(22, 5).Value = 1.046845016649588
$ws.Cells.Item(22, 6).Value = 1.054503653239638
$ws.Cells.Item(22, 9).Value = 1.033397503852102
$ws.Cells.Item(22, 10).Value = 1.046523330715005
$ws.Cells.Item(22, 11).Value = 1.042419669330028
$ws.Cells.Item(22, 12).Value = 1.050512703232796
$ws.Cells.Item(22, 13).Value = 1.058142649199402

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.040324241398519
$ws.Cells.Item(23, 4).Value = 1.039154336167626
$ws.Cells.Item(23, 5).Value = 1.04739664448864
$ws.Cells.Item(23, 6).Value = 1.055124648839804
$ws.Cells.Item(23, 9).Value = 1.03351643125809
$ws.Cells.Item(23, 10).Value = 1.046993146119551
$ws.Cells.Item(23, 11).Value = 1.042783255762362
$ws.Cells.Item(23, 12).Value = 1.050994814552325
$ws.Cells.Item(23, 13).Value = 1.058694470406689

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.042692503622208
$ws.Cells.Item(24, 4).Value = 1.040855752281336
$ws.Cells.Item(24, 5).Value = 1.049564654669631
$ws.Cells.Item(24, 6).Value = 1.057566147817943
$ws.Cells.Item(24, 9).Value = 1.03397976094024
$ws.Cells.Item(24, 10).Value = 1.048837382896149
$ws.Cells.Item(24, 11).Value = 1.044208955779466
$ws.Cells.Item(24, 12).Value = 1.052888052708425
$ws.Cells.Item(24, 13).Value = 1.060862636621174

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.045431251036328
$ws.Cells.Item(25, 4).Value = 1.042821198021078
$ws.Cells.Item(25, 5).Value = 1.05207320139463
$ws.Cells.Item(25, 6).Value = 1.060392930192735
$ws.Cells.Item(25, 9).Value = 1.034507420008806
$ws.Cells.Item(25, 10).Value = 1.050966686135081
$ws.Cells.Item(25, 11).Value = 1.045851823938286
$ws.Cells.Item(25, 12).Value = 1.055075451771201
$ws.Cells.Item(25, 13).Value = 1.063370125040097
